$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codespace")

# Add new row 5: a mapping entry for the typo "ESPG" -> "EPSG"
$ws.Range("A5").Value = "ESPG"
$ws.Range("B5").Value = "EPSG"

# Update the saved selection to A6 (next empty row), matching the diff
$ws.Range("A6").Select()
